$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.212.97'
$ws.Range('E2').Value = '  +0.06%  '

$ws.Range('D3').Value = '2.849.62'
$ws.Range('E3').Value = '  +2.10%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '361.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.94%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '113.92'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.60%  '

$ws.Range('E7').Value = '  +3.29%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.07%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.603'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.02%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.68'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.01%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0862'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.82%  '

$ws.Range('E12').Value = '  +1.27%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.05'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.51%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.78'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.13%  '

$ws.Range('D15').Value = '3.300.40'
$ws.Range('E15').Value = '  +2.19%  '

$ws.Range('D16').Value = '2.846.65'
$ws.Range('E16').Value = '  +1.51%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.906'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.25%  '

$ws.Range('D18').Value = '52.033.18'
$ws.Range('E18').Value = '  +0.00%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.72%  '

$ws.Range('E20').Value = '  -1.81%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.60'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.88%  '

$ws.Range('D22').Value = '0.0₃0994'
$ws.Range('E22').Value = '  +0.73%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.39'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.08%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '268.75'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.39%  '

$ws.Range('E25').Value = '  +0.59%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.22'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.27%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.07%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.44'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.53%  '

$ws.Range('E29').Value = '  +1.56%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '53.67'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.23%  '

$ws.Range('E31').Value = '  -1.80%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0470'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +25.57%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.16'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.14%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.89'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.10%  '

$ws.Range('E35').Value = '  +8.33%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0847'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.53%  '

$ws.Range('E37').Value = '  +0.07%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.29'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.81%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.08'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.08%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.38'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.21%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '23.85'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.22%  '

$ws.Range('E42').Value = '  +1.12%  '

$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '128.34'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.99%  '

$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.57'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.23%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.27'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.42%  '

$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '2.114.30'
$ws.Range('E46').Value = '  +0.43%  '

$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.40'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.85%  '

$ws.Range('E49').Value = '  +10.74%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.87'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.90%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.07'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.34%  '
